# Generate Report for Handback
#
# For both locale sheets (zh-cn, de-de), this:
#   - updates the Status column (C) text from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two data rows
#   - fills in the "Latest Target File" (F) and "Latest Handback File" (G)
#     columns with hyperlinked filenames (mirroring the existing
#     "Latest Handoff File" entries in column D, since the handback
#     round-trips the same files)
#   - stamps the "Latest Handback DateTime" (H) column with the real
#     handback timestamp instead of the epoch placeholder

$wb = $excel.ActiveWorkbook

$hyperlinkFont = 15570276  # BGR for OOXML rgb FF6495ED (matches existing hyperlink style)

function Set-HandbackRow {
    param($ws, $row, $mdName, $mdUrl, $xlfName, $xlfUrl, $handbackDateTime)

    # Status
    $ws.Range("C$row").Value2 = "Handed back: in sync with en-US"

    # Latest Target File (F) -> same source .md as the handoff file
    $fCell = $ws.Range("F$row")
    $fCell.Value2 = $mdName
    $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdName) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkFont

    # Latest Handback File (G) -> the translated .xlf that came back
    $gCell = $ws.Range("G$row")
    $gCell.Value2 = $xlfName
    $ws.Hyperlinks.Add($gCell, $xlfUrl, "", "", $xlfName) | Out-Null
    $gCell.Font.Underline = 2
    $gCell.Font.Color = $hyperlinkFont

    # Latest Handback DateTime (H)
    $ws.Range("H$row").Value2 = $handbackDateTime
}

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $zh 2 `
    "04678046-4acf-41ea-ba88-41ade6e9999c.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a7914207129647b570263f54f7a45d9d535f1882/e2e/04678046-4acf-41ea-ba88-41ade6e9999c.md" `
    "04678046-4acf-41ea-ba88-41ade6e9999c.d7445fdfeb80f7bddfe3c27c7b21de424f9046df.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49ccdb2df0f7782aa0fa7a1ae92ae52a0eee17db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/04678046-4acf-41ea-ba88-41ade6e9999c.d7445fdfeb80f7bddfe3c27c7b21de424f9046df.zh-cn.xlf" `
    "2016-03-22 19:05:40"

Set-HandbackRow $zh 3 `
    "44e5c74c-4e04-48d2-a31f-602ffd7cd41c.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a7914207129647b570263f54f7a45d9d535f1882/e2e/44e5c74c-4e04-48d2-a31f-602ffd7cd41c.md" `
    "44e5c74c-4e04-48d2-a31f-602ffd7cd41c.3ffa520b476fb3eddf8b12577ef49f603accd16f.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49ccdb2df0f7782aa0fa7a1ae92ae52a0eee17db/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/44e5c74c-4e04-48d2-a31f-602ffd7cd41c.3ffa520b476fb3eddf8b12577ef49f603accd16f.zh-cn.xlf" `
    "2016-03-22 19:05:40"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

Set-HandbackRow $de 2 `
    "04678046-4acf-41ea-ba88-41ade6e9999c.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a7914207129647b570263f54f7a45d9d535f1882/e2e/04678046-4acf-41ea-ba88-41ade6e9999c.md" `
    "04678046-4acf-41ea-ba88-41ade6e9999c.d7445fdfeb80f7bddfe3c27c7b21de424f9046df.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cde936afe080c75ec5f7c49549d770c99d4d9d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/04678046-4acf-41ea-ba88-41ade6e9999c.d7445fdfeb80f7bddfe3c27c7b21de424f9046df.de-de.xlf" `
    "2016-03-22 19:05:47"

Set-HandbackRow $de 3 `
    "44e5c74c-4e04-48d2-a31f-602ffd7cd41c.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a7914207129647b570263f54f7a45d9d535f1882/e2e/44e5c74c-4e04-48d2-a31f-602ffd7cd41c.md" `
    "44e5c74c-4e04-48d2-a31f-602ffd7cd41c.3ffa520b476fb3eddf8b12577ef49f603accd16f.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cde936afe080c75ec5f7c49549d770c99d4d9d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/44e5c74c-4e04-48d2-a31f-602ffd7cd41c.3ffa520b476fb3eddf8b12577ef49f603accd16f.de-de.xlf" `
    "2016-03-22 19:05:47"

Write-Host "Handback report generated."
